# Generate Report for Handback
#
# The localization handback just completed for the zh-cn and de-de
# targets. Refresh the status report: flip the overall status from
# "Ready for handoff" to "Handed back: in sync with en-US", stamp the
# handback timestamps, link the produced target files, and widen the
# columns that now hold longer content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status text for both locales -------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# --- zh-cn sheet ------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# Latest Target File (I) now links to the source .md that was handed back
$zhcn.Range("I2").Value = "674e4054-463e-4f57-b03f-a75f83fdb0ff.md"
$zhcn.Range("I2").Font.Underline = $true
$zhcn.Range("I3").Value = "c57c4e29-70fe-47d4-9162-3e3b2459bc17.md"
$zhcn.Range("I3").Font.Underline = $true

# Latest Handback File (J) now references the generated xliff
$zhcn.Range("J2").Value = "674e4054-463e-4f57-b03f-a75f83fdb0ff.fc4e56ee4a7d451e7532a349b2707d279729f771.zh-cn.xlf"
$zhcn.Range("J3").Value = "c57c4e29-70fe-47d4-9162-3e3b2459bc17.e50ac1374d38381e5ab2965ff4581a0975325269.zh-cn.xlf"

# Latest Handback DateTime (K) - handback completed
$zhcn.Range("K2").Value = "2016-09-01 13:09:25"
$zhcn.Range("K3").Value = "2016-09-01 13:09:25"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3350a9e6a80a850e56b8358d6380cfd04b702a8a/e2e/674e4054-463e-4f57-b03f-a75f83fdb0ff.md", "", "", "674e4054-463e-4f57-b03f-a75f83fdb0ff.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3350a9e6a80a850e56b8358d6380cfd04b702a8a/e2e/c57c4e29-70fe-47d4-9162-3e3b2459bc17.md", "", "", "c57c4e29-70fe-47d4-9162-3e3b2459bc17.md") | Out-Null

$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# --- de-de sheet --------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("I2").Value = "674e4054-463e-4f57-b03f-a75f83fdb0ff.md"
$dede.Range("I2").Font.Underline = $true
$dede.Range("I3").Value = "c57c4e29-70fe-47d4-9162-3e3b2459bc17.md"
$dede.Range("I3").Font.Underline = $true

$dede.Range("J2").Value = "674e4054-463e-4f57-b03f-a75f83fdb0ff.fc4e56ee4a7d451e7532a349b2707d279729f771.de-de.xlf"
$dede.Range("J3").Value = "c57c4e29-70fe-47d4-9162-3e3b2459bc17.e50ac1374d38381e5ab2965ff4581a0975325269.de-de.xlf"

$dede.Range("K2").Value = "2016-09-01 13:09:32"
$dede.Range("K3").Value = "2016-09-01 13:09:32"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3350a9e6a80a850e56b8358d6380cfd04b702a8a/e2e/674e4054-463e-4f57-b03f-a75f83fdb0ff.md", "", "", "674e4054-463e-4f57-b03f-a75f83fdb0ff.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3350a9e6a80a850e56b8358d6380cfd04b702a8a/e2e/c57c4e29-70fe-47d4-9162-3e3b2459bc17.md", "", "", "c57c4e29-70fe-47d4-9162-3e3b2459bc17.md") | Out-Null

$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17

Write-Host "Handback report regenerated."
